# Enhance logging system configuration: append the newest log record (row 77)
# to each of the four MID_* sheets, mirroring the formatting of row 76.

$wb = $excel.ActiveWorkbook

$rowsData = @{
    1 = @{
        A = 45863.46560185185
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x3C"
        E = "0x07"
        F = 400
        G = 568631262647113031352320.0
        H = 316
        I = 7
    }
    2 = @{
        A = 45863.46560185185
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x38"
        E = "0x19"
        F = 380
        G = 568432987514711010443264.0
        H = 312
        I = 25
    }
    3 = @{
        A = 45863.46560185185
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x61"
        E = "0x15"
        F = 110
        G = 568631262647113031352320.0
        H = 97
        I = 15
    }
    4 = @{
        A = 45863.46560185185
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x77"
        E = "0x9"
        F = 130
        G = 568631262647113031352320.0
        H = 119
        I = 9
    }
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $data = $rowsData[$sheetIndex]

    $newRow = 77

    # Column A: timestamp, keep same number format as the row above (row 76)
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    # Columns B-E: hex-byte strings stored as text
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: numeric values
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
